$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reduce row 5 readings to 2-decimal "custom accuracy" ---
$ws.Range("B5").Value = 17.17
$ws.Range("C5").Value = 12.79
$ws.Range("D5").Value = 1.12
$ws.Range("E5").Value = 37.61
$ws.Range("G5").Value = 13.46
$ws.Range("H5").Value = 49.73
$ws.Range("I5").Value = 20.94
$ws.Range("J5").Value = 9.23
$ws.Range("K5").Value = 13.59
$ws.Range("M5").Value = 16.08
$ws.Range("N5").Value = 4.24
$ws.Range("O5").Value = 13.53
$ws.Range("P5").Value = 19.19
$ws.Range("R5").Value = 0.75
$ws.Range("S5").Value = 0.73
$ws.Range("T5").Value = 198.68
$ws.Range("U5").Value = 37.77
$ws.Range("V5").Value = 12.49
$ws.Range("Y5").Value = 2.03
$ws.Range("Z5").Value = 24.72
$ws.Range("AA5").Value = 11.03
$ws.Range("AB5").Value = 9.83
$ws.Range("AC5").Value = 11.54
$ws.Range("AD5").Value = 15.85
$ws.Range("AG5").Value = 6.98
$ws.Range("AH5").Value = 15.62

# --- 2. Drop the now-redundant row 6 (data trimmed to 1000 rows upstream) ---
$ws.Rows.Item(6).Delete()

# --- 3. Narrow a handful of data columns by one character ---
$narrowCols = @(3, 7, 11, 15, 22, 27, 29)
foreach ($col in $narrowCols) {
    $ws.Columns.Item($col).ColumnWidth = 6.166666667
}
